# Update "想去人数" (want-to-go count) figures refreshed by the gh-pages
# data scrape (commit 456a3b4).
#
# Sheet "展览"  (sheet1): rows 2,5,6,8,9,10
# Sheet "全部类型" (sheet4): rows 2,5,6,8,9,11 (same events, offset by the
#   two extra rows present only on the "全部类型" sheet)

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetExhibit.Range("F2").Value = 318
$sheetExhibit.Range("F5").Value = 4672
$sheetExhibit.Range("F6").Value = 364
$sheetExhibit.Range("F8").Value = 287
$sheetExhibit.Range("F9").Value = 723
$sheetExhibit.Range("F10").Value = 204

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 318
$sheetAll.Range("F5").Value = 4672
$sheetAll.Range("F6").Value = 364
$sheetAll.Range("F8").Value = 287
$sheetAll.Range("F9").Value = 723
$sheetAll.Range("F11").Value = 204
